$wb = $excel.ActiveWorkbook

# hunk 0: sheet ALC row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 12000
$ws.Range("J21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("N21").Value = -12936

# hunk 1: sheet ALC row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 12000
$ws.Range("J23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("N23").Value = -12468

# hunk 2: sheet ALC row 26
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J26").Value = 50000
$ws.Range("L26").Value = 50000
$ws.Range("N26").Value = -50688

# hunk 3: sheet ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3585.7144
$ws.Range("I74").Value = 3516.6667
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3516.6667
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2580.6667
$ws.Range("N74").Value = -5872

# hunk 4: sheet ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3585.7144
$ws.Range("I77").Value = 3516.6667
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 17583.3335
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -12903.3335
$ws.Range("N77").Value = -29360

# hunk 5: sheet ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1070.8379
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 1086.6945
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 3260.0835
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -13260.0835

# hunk 6: sheet ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 38308
$ws.Range("J133").Value = 38308
$ws.Range("L133").Value = 38308
$ws.Range("N133").Value = -48428

# hunk 7: sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1244.75
$ws.Range("I137").Value = 1170.8334
$ws.Range("K137").Value = 3512.5002
$ws.Range("M137").Value = -962.5001999999999

# hunk 8: sheet ARM row 21
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 23252.5
$ws.Range("I21").Value = 23252.5
$ws.Range("K21").Value = 23252.5
$ws.Range("M21").Value = -22878.5

# hunk 9: sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24060.451
$ws.Range("I32").Value = 3834.7708
$ws.Range("K32").Value = 3834.7708
$ws.Range("M32").Value = -3547.7708

# hunk 10: sheet ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2050
$ws.Range("I110").Value = 1877.75
$ws.Range("K110").Value = 1877.75
$ws.Range("M110").Value = 167.25

# hunk 11: sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3988.2
$ws.Range("I122").Value = 1980.3334
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 5941.0002
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -3491.0002
$ws.Range("N122").Value = -25900

# hunk 12: sheet ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 63615.25
$ws.Range("J133").Value = 63615.25
$ws.Range("L133").Value = 63615.25
$ws.Range("N133").Value = -68675.25

# hunk 13: sheet ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 48328.89
$ws.Range("J139").Value = 48328.89
$ws.Range("L139").Value = 48328.89
$ws.Range("N139").Value = -58608.89

# hunk 14: sheet ARM row 141
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360

# hunk 15: sheet BSM row 133
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 38780
$ws.Range("J133").Value = 38780
$ws.Range("L133").Value = 38780
$ws.Range("N133").Value = -48900

# hunk 16: sheet CRP row 15
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

# hunk 17: sheet CRP row 23
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 47501.75
$ws.Range("I23").Value = 45004.5
$ws.Range("K23").Value = 45004.5
$ws.Range("M23").Value = -44764.5

# hunk 18: sheet CRP row 27
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 47501.75
$ws.Range("I27").Value = 45004.5
$ws.Range("K27").Value = 45004.5
$ws.Range("M27").Value = -44812.5

# hunk 19: sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = $null

# hunk 20: sheet CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 71380
$ws.Range("J140").Value = 71380
$ws.Range("L140").Value = 71380
$ws.Range("N140").Value = -81740

# hunk 21: sheet CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 58823690
$ws.Range("I2").Value = 26.727272
$ws.Range("J2").Value = 166667070
$ws.Range("K2").Value = 160.363632
$ws.Range("L2").Value = 1000002420
$ws.Range("M2").Value = -47.363632
$ws.Range("N2").Value = -1000002646

# hunk 22: sheet CUL row 49
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 1310.1666
$ws.Range("I49").Value = 728.6667
$ws.Range("K49").Value = 2186.0001
$ws.Range("M49").Value = -2030.0001

# hunk 23: sheet CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 9000
$ws.Range("J58").Value = 9000
$ws.Range("L58").Value = 27000
$ws.Range("N58").Value = -27256

# hunk 24: sheet CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 17619.143
$ws.Range("I68").Value = 30532.75
$ws.Range("K68").Value = 91598.25
$ws.Range("M68").Value = -90787.25

# hunk 25: sheet CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 17619.143
$ws.Range("I71").Value = 30532.75
$ws.Range("K71").Value = 274794.75
$ws.Range("M71").Value = -270738.75

# hunk 26: sheet GSM row 20
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = $null

# hunk 27: sheet GSM row 27
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null

# hunk 28: sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2438.7896
$ws.Range("I122").Value = 2255.1177
$ws.Range("K122").Value = 6765.353099999999
$ws.Range("M122").Value = -4315.353099999999

# hunk 29: sheet GSM row 137
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 48000
$ws.Range("J137").Value = 48000
$ws.Range("L137").Value = 48000
$ws.Range("N137").Value = -58200

# hunk 30: sheet GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 75741
$ws.Range("J138").Value = 75741
$ws.Range("L138").Value = 75741
$ws.Range("N138").Value = -86021

# hunk 31: sheet GSM row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 28800
$ws.Range("J139").Value = 28800
$ws.Range("L139").Value = 28800
$ws.Range("N139").Value = -39080

# hunk 32: sheet LTW row 14
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 26502
$ws.Range("I14").Value = 26502
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 26502
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -26330
$ws.Range("N14").Value = $null

# hunk 33: sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 17650
$ws.Range("I22").Value = 966.6667
$ws.Range("K22").Value = 966.6667
$ws.Range("M22").Value = -671.6667

# hunk 34: sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 17650
$ws.Range("I27").Value = 966.6667
$ws.Range("K27").Value = 966.6667
$ws.Range("M27").Value = -859.6667

# hunk 35: sheet LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 900
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 833.3333
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 833.3333
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1209.3333

# hunk 36: sheet LTW row 92
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null

# hunk 37: sheet WVR row 5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1000497.5
$ws.Range("I5").Value = 995
$ws.Range("K5").Value = 995
$ws.Range("M5").Value = -883

# hunk 38: sheet WVR row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = -10480

# hunk 39: sheet WVR row 26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = $null

# hunk 40: sheet WVR row 32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 500
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -183
$ws.Range("N32").Value = $null

# hunk 41: sheet WVR row 40
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = $null

# hunk 42: sheet WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3086.7222
$ws.Range("I81").Value = 1887.3636
$ws.Range("K81").Value = 3774.7272
$ws.Range("M81").Value = -2713.7272

# hunk 43: sheet WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3086.7222
$ws.Range("I84").Value = 1887.3636
$ws.Range("K84").Value = 18873.636
$ws.Range("M84").Value = -13569.636

# hunk 44: sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 144429
$ws.Range("I122").Value = 201360.6
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 604081.8
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -601631.8
$ws.Range("N122").Value = -11200

# hunk 45: sheet WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 92089.91
$ws.Range("I126").Value = 167448.17
$ws.Range("J126").Value = 1660
$ws.Range("K126").Value = 502344.51
$ws.Range("L126").Value = 4980
$ws.Range("M126").Value = -499874.51
$ws.Range("N126").Value = -9920
